$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'mens compression tights white',
    'mens compression workout pants',
    'mens crossfit pants',
    'mens cycling pants long for cold weather',
    'mens fitness wear',
    'mens guide pro pants',
    'mens high waist trousers',
    'mens hiking leggings',
    'mens leggings compression nike',
    'mens leggings compression winter',
    'mens leggings for winter',
    'mens leggings xs',
    'mens light blue compression pants',
    'mens light blue leggings',
    'mens lycra running pants',
    'mens medical compression leggings',
    'mens muscle pants',
    'mens neoprene',
    'mens neoprene pants',
    'mens red compression tights',
    'mens running compression pants',
    'mens running compression tights',
    'mens running knee compression',
    'mens running leggings winter',
    'mens running long pants',
    'mens running tights',
    'mens running tights 3 4',
    'mens running tights capri',
    'mens running tights cold weather',
    'mens running tights long',
    'mens running tights pocket',
    'mens running tights with pocket',
    'mens running underwear cold',
    'mens see through leggings',
    'mens ski wear',
    'mens skin tight leggings',
    'mens skins tights',
    'mens soccer pants',
    'mens soccer training pants',
    'mens support underwear',
    'mens tear away athletic pants',
    'mens tight',
    'mens tights basketball',
    'mens tights capri',
    'mens tights cold weather',
    'mens tights leggings',
    'mens tights thermal',
    'mens tights winter thermal',
    'mens tights with feet',
    'mens tights with knee pads',
    'mens track leggings',
    'mens training tights',
    'mens warm leggings compression',
    'mens winter compression pants',
    'mens winter workout gear',
    'mens workout back brace',
    'mens workout tights',
    'mens workout tights pants',
    'mesh active leggings',
    'mesh tights women',
    'mid calf leggings',
    'mid calf leggings for women',
    'mid calf workout leggings for women',
    'mid high',
    'mid length',
    'mid length shorts women',
    'mid rise leggings',
    'mid rise leggings for women',
    'mid rise tights',
    'mid rise workout leggings',
    'mid skirt winter',
    'mid thigh athletic shorts for women',
    'mid thigh bike shorts',
    'mid thigh black shorts for women',
    'mid thigh skirt',
    'mid thigh skirts for women',
    'mid tops for women',
    'mid trail runners',
    'mid waist leggings for women',
    'midrise hiking pants',
    'midrise lift',
    'mids and highs amplifier',
    'midsection compression',
    'mil shake',
    'mil tech',
    'mild cigarettes',
    'mile high media',
    'mile tracker for bike',
    'mile tracker for running',
    'mile tracking',
    'mind design wash cloths',
    'mind rise',
    'mix martial arts equipment',
    'mma compression shorts',
    'mobility tape',
    'modal tights',
    'modal underwear women high',
    'model woman',
    'modest gym clothes for women',
    'modest leggings for women'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
